$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells are stored as literal text (they use "." as a
# thousands separator in some rows, e.g. "61.834.22"), never as numbers.
# Force text typing via NumberFormat "@" while writing, then restore the
# cell style to Normal/General so no stray numeric/text style sticks around.
function Set-DText($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-DText "D2" "61.834.22"
$ws.Range("E2").Value = "  +1.25%  "
# Row 3
Set-DText "D3" "3.459.92"
$ws.Range("E3").Value = "  +2.98%  "
# Row 4
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
Set-DText "D5" "580.08"
$ws.Range("E5").Value = "  +1.64%  "
# Row 6
Set-DText "D6" "148.31"
$ws.Range("E6").Value = "  +9.21%  "
# Row 7
Set-DText "D7" "3.460.94"
$ws.Range("E7").Value = "  +3.08%  "
# Row 8
$ws.Range("E8").Value = "  +0.00%  "
# Row 9
$ws.Range("E9").Value = "  +1.23%  "
# Row 10
Set-DText "D10" "7.77"
$ws.Range("E10").Value = "  +4.02%  "
# Row 11
Set-DText "D11" "0.124"
$ws.Range("E11").Value = "  +1.48%  "
# Row 12
$ws.Range("E12").Value = "  +1.67%  "
# Row 13
Set-DText "D13" "4.051.47"
$ws.Range("E13").Value = "  +3.00%  "
# Row 14
Set-DText "D14" "28.10"
$ws.Range("E14").Value = "  +8.33%  "
# Row 15
$ws.Range("E15").Value = "  -0.22%  "
# Row 16
Set-DText "D16" "0.0000175"
$ws.Range("E16").Value = "  +1.88%  "
# Row 17
Set-DText "D17" "3.451.91"
$ws.Range("E17").Value = "  +2.67%  "
# Row 18
Set-DText "D18" "61.858.33"
$ws.Range("E18").Value = "  +1.10%  "
# Row 19
Set-DText "D19" "6.34"
$ws.Range("E19").Value = "  +9.22%  "
# Row 20
Set-DText "D20" "14.38"
$ws.Range("E20").Value = "  +2.96%  "
# Row 21
Set-DText "D21" "9.47"
$ws.Range("E21").Value = "  +2.52%  "
# Row 22
Set-DText "D22" "385.37"
$ws.Range("E22").Value = "  +2.27%  "
# Row 23
$ws.Range("E23").Value = "  +3.22%  "
# Row 24
Set-DText "D24" "3.596.53"
$ws.Range("E24").Value = "  +3.01%  "
# Row 25
$ws.Range("E25").Value = "  +0.18%  "
# Row 26
$ws.Range("E26").Value = "  +0.97%  "
# Row 27
Set-DText "D27" "72.54"
$ws.Range("E27").Value = "  +2.12%  "
# Row 28
$ws.Range("E28").Value = "  -1.36%  "
# Row 29
$ws.Range("E29").Value = "  +9.45%  "
# Row 30
Set-DText "D30" "7.84"
$ws.Range("E30").Value = "  +5.84%  "
# Row 31
Set-DText "D31" "1.54"
$ws.Range("E31").Value = "  -12.90%  "
# Row 32
$ws.Range("E32").Value = "  -0.59%  "
# Row 33
Set-DText "D33" "8.25"
$ws.Range("E33").Value = "  +1.50%  "
# Row 34
$ws.Range("E34").Value = "  +1.99%  "
# Row 35
$ws.Range("E35").Value = "  +0.01%  "
# Row 36
$ws.Range("E36").Value = "  +2.07%  "
# Row 37
$ws.Range("E37").Value = "  +4.55%  "
# Row 38
Set-DText "D38" "5.20"
$ws.Range("E38").Value = "  +0.45%  "
# Row 39
$ws.Range("E39").Value = "  +2.71%  "
# Row 40
Set-DText "D40" "167.02"
$ws.Range("E40").Value = "  +1.18%  "
# Row 41
Set-DText "D41" "0.0789"
$ws.Range("E41").Value = "  +4.67%  "
# Row 42
Set-DText "D42" "0.798"
$ws.Range("E42").Value = "  +3.89%  "
# Row 43
Set-DText "D43" "26.00"
$ws.Range("E43").Value = "  +9.81%  "
# Row 44
$ws.Range("E44").Value = "  +2.05%  "
# Row 45
$ws.Range("E45").Value = "  -0.11%  "
# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-DText "D46" "4.50"
$ws.Range("E46").Value = "  +2.72%  "
# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-DText "D47" "42.36"
$ws.Range("E47").Value = "  +2.26%  "
# Row 48
$ws.Range("E48").Value = "  -2.19%  "
# Row 49
Set-DText "D49" "2.600.23"
$ws.Range("E49").Value = "  +10.98%  "
# Row 50
Set-DText "D50" "6.95"
$ws.Range("E50").Value = "  +2.36%  "
# Row 51
Set-DText "D51" "23.33"
$ws.Range("E51").Value = "  +0.67%  "
